$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: db-update task finished later than planned; update actual
#     hours, actual delivery date and the comment describing the extra work
#     (toplist support from Rakib). ---
$ws.Range("P19").Value = 63
$ws.Range("Q19").Value = "2020-03-12 Pending"
$ws.Range("T19").Value = "Greate trouble making it work to update db via react app. Then great trouble sending all user fields to db. Finally got through with support from Rakib"

# --- Row 20: fill in actual hours / actual delivery date / weekday and the
#     plan-vs-act formula (previously still blank/pending). ---
$ws.Range("P20").Value = 3
$ws.Range("Q20").Value = 43902
$ws.Range("R20").Value = "Friday"
$ws.Range("S20").Formula = "=M20-P20"

# --- Row 21: same for the following task. ---
$ws.Range("P21").Value = 28
$ws.Range("Q21").Value = "2020-03-23 Pending"
$ws.Range("R21").Value = "Monday"
$ws.Range("S21").Formula = "=M21-P21"

# --- Row 29 "Topplistesida": added the toplist endpoint to the api, so its
#     actual hours / actual delivery date / weekday can now be filled in.
#     Q29 had no date formatting yet (unlike Q20/Q21), so copy the date
#     number-format from the cell above before writing the value. ---
$ws.Range("P29").Value = 8
$ws.Range("Q20").Copy()
$ws.Range("Q29").PasteSpecial(-4122)
$ws.Range("Q29").Value = 43921
$ws.Range("R29").Value = "Tuesday"

# Keep the on-screen selection in sync with where the edits ended up.
$ws.Range("P29").Select()
